# Update the "Förändrad" (changed) date column C for every data row
# (rows 2-236) from 2023-09-15 (45184) to 2023-09-17 (45186), and add a
# friendly-name second argument (the "Beteckning" in column A) to every
# HYPERLINK() formula found in columns S, T, V, W, X, Y (rows 2-21, the
# only rows that currently carry those link formulas).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow  = 236
$linkCols     = @("S", "T", "V", "W", "X", "Y")

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {

    # Column C: bump the "Förändrad" serial date from 45184 (2023-09-15)
    # to 45186 (2023-09-17) — only touch rows that actually carried the
    # old value (defensive; all data rows do in this workbook).
    $cCell = $ws.Range("C$r")
    if ($cCell.Text -eq "2023-09-15") {
        $cCell.Value = 45186
    }

    # Columns S/T/V/W/X/Y: append the Beteckning (column A) as the
    # second, friendly-name argument of the HYPERLINK() formula, when
    # such a formula is present on this row.
    $beteckning = $ws.Range("A$r").Text

    foreach ($col in $linkCols) {
        $cell = $ws.Range("$col$r")
        $f = $cell.Formula

        if ($f -and $f.Length -gt 0 -and $f.ToUpper().Contains("HYPERLINK(") -and -not $f.Contains(",")) {
            $trimmed = $f.TrimEnd()
            $newFormula = $trimmed.Substring(0, $trimmed.Length - 1) + ', "' + $beteckning + '")'
            $cell.Formula = $newFormula
        }
    }
}
